$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.519.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.596.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.62"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.346"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.051.75"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "60.562.87"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.77"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.606.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "353.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.717.53"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.930"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0840"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.26"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +10.42%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.58"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.24"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.942"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +11.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.19"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.35"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.77"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.51"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.976.12"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.90%  "
